# Update "paises.xlsx" ("Pais" sheet) with refreshed COVID-19 country data
# and updated timestamp, per commit "Update countries & provincias Spain".
#
# The source data is sorted by column B ("Casos totales") descending. Refreshed
# counts for a handful of countries change their relative ranking, which is why
# some adjacent rows end up with different country names after the refresh
# (the row's data travels with the country, not with the row number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write one data row (columns B..H) in one shot ---------------
function Set-Row {
    param($row, $pais, $b, $c, $d, $e, $f, $g, $h)
    $ws.Cells.Item($row, 1).Value = $pais
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- A1: refreshed "last updated" timestamp -------------------------------
$ws.Range("A1").Value = "Datos actualizados a 17 de Junio de 2020 a las 07:16"

# --- Plain data refreshes (country keeps its row) -------------------------
Set-Row 5   "Brasil"      928834 0   477364 406014 0 0 45456
Set-Row 67  "Honduras"    9656   478 1075   8251   0 8 330
Set-Row 76  "Uzbekistan"  5561   68  4096   1446   0 0 19
Set-Row 93  "Tailandia"   3135   0   2996   81     0 0 58
Set-Row 96  "Kirguistan"  2562   90  1902   630    0 1 30
Set-Row 158 "Vietnam"     335    1   325    10     0 0 0
Set-Row 184 "Butan"       67     0   24     43     0 0 0
Set-Row 198 "Belice"      22     1   16     4      0 0 2

# --- Haiti jumps ahead of Republica de Yibuti with fresh numbers ----------
Set-Row 81  "Haiti"                4547 106 24   4443 0 4 80
Set-Row 82  "Republica de Yibuti"  4539 0   3324 1172 0 0 43

# --- Tied-total pairs whose relative order flips in the refreshed feed ----
Set-Row 206 "Groenlandia"                13 0 13 0 0 0 0
Set-Row 207 "Islas Malvinas"             13 0 13 0 0 0 0

Set-Row 210 "Seychelles"                 11 0 11 0 0 0 0
Set-Row 211 "Montserrat"                 11 0 10 0 0 0 1

Set-Row 213 "Papua Nueva Guinea"         8  0 8  0 0 0 0
Set-Row 214 "Islas Virgenes Britanicas"  8  0 7  0 0 0 1
